$wb = $excel.ActiveWorkbook

# --- "Trends Status" sheet ---
$wsTrends = $wb.Worksheets.Item("Trends Status")
$wsTrends.Range("C3").Value = 2          # Decline / Current species (no.)
$wsTrends.Range("E3").Value = 66.7       # Decline / Current species conclusive (perc.)
$wsTrends.Range("E4").Value = 33.3       # Stable / Current species conclusive (perc.)
$wsTrends.Range("C7").Value = 12         # Trend Inconclusive / Current species (no.)

# --- "Species qualification" sheet ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("C4").Value = 3            # Current Analysis / With conclusive trends

# --- "Interannual update - High Pri" sheet ---
# Insert a new row 3 ("Trend Different") before the existing "IUCN" row,
# which pushes IUCN down to row 4 and updates its counts.
$wsInter = $wb.Worksheets.Item("Interannual update - High Pri")
$wsInter.Rows("3").Insert()

$wsInter.Range("A3").Value = "Trend Different"
$wsInter.Range("B3").Value = 1
$wsInter.Range("C3").Value = 1

$wsInter.Range("B4").Value = 20
$wsInter.Range("C4").Value = 19.4
